$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 changes
$ws.Range("M3").Value = 1.17
$ws.Range("N3").Value = 5

# Row 6 changes
$ws.Range("G6").Value = 1.5
$ws.Range("H6").Value = 4
$ws.Range("I6").Value = 6.25
$ws.Range("J6").Value = 2.05
$ws.Range("K6").Value = 2.4
$ws.Range("U6").Value = 1.8
$ws.Range("V6").Value = 1.91
$ws.Range("AG6").Value = 19
$ws.Range("AH6").Value = 34
$ws.Range("AJ6").Value = 67
$ws.Range("AQ6").Value = 21
$ws.Range("AS6").Value = 126
$ws.Range("AW6").Value = 7.5

$wb.Save()
